$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"8.43769498715119e-17"
$ws.Range("F2").Value = [double]"1.576516694967722e-16"

$ws.Range("E3").Value = [double]"0.01149607432881209"
$ws.Range("F3").Value = [double]"0.01543341287607883"

$ws.Range("E4").Value = [double]"-0.006754119518834534"
$ws.Range("F4").Value = [double]"-0.01028894191738561"

$ws.Range("E5").Value = [double]"-0.005058824988670463"
$ws.Range("F5").Value = [double]"0.04115576766954331"

$ws.Range("E6").Value = [double]"0.001108556945570738"
$ws.Range("F6").Value = [double]"0.04115576766954331"

$ws.Range("E7").Value = [double]"-0.0007916867668774308"
$ws.Range("F7").Value = [double]"-0.08745600629777898"
